# Release: Increment splash screen / about box / plug-in version to 1.0.0,
# and refresh the cached "datetimeFigureOut" date field text (8/3/2010 -> 8/17/2010)
# on the slide master and every slide layout's Date Placeholder.

$p = $ppt.ActivePresentation

$oldDate = "8/3/2010"
$newDate = "8/17/2010"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's Date Placeholder.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1: bump the "Release 0.3.0" text box to "Release 1.0.0" without
# touching the separately-formatted "Release " run.
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $idx = $tr.Text.IndexOf("0.3.0")
        if ($idx -ge 0) {
            # Remember the autofit-derived box size; re-running the text
            # measurement after the (same-length) substitution can nudge
            # the cached height by a few EMUs, which the real edit did not
            # do, so put it back exactly afterwards.
            $origHeight = $sh.Height
            $origWidth = $sh.Width
            $sub = $tr.Characters($idx + 1, 5)
            $sub.Text = "1.0.0"
            $sh.Height = $origHeight
            $sh.Width = $origWidth
        }
    }
}
